# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — style matches the other header cells (bold, bordered,
# centered, top-aligned), same as cell AC1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record values for every player row (2 through 45).
$wins = 78
$losses = 84
$ties = 0

for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
